$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G: "Material Type" header + "DNA:Genomic" for every data row
$ws.Range("G1").Value = "Material Type"
$ws.Range("G2:G24").Value = "DNA:Genomic"

# Give the new header cell the same bold/centered/filled look as the other
# header cells (copy F1's resolved style), then recolor it to a black fill
# with bold white text.
$ws.Range("F1").Copy()
$headerCell = $ws.Range("G1")
$headerCell.PasteSpecial(-4122)
$headerCell.Font.Color = 16777215
$headerCell.Interior.Color = 0
$excel.CutCopyMode = $false

# Body cells G2:G24 get the same centered style already used in column C
$bodyRange = $ws.Range("G2:G24")
$bodyRange.HorizontalAlignment = -4108

# Selection matches the authored state
$ws.Range("G1:G24").Select()
